$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.010.71'
$ws.Range("E2").Value = '  +1.40%  '
$ws.Range("D3").Value = '1.642.58'
$ws.Range("E3").Value = '  +0.58%  '
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '213.54'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.66%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.525'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.28%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.999'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.06%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '23.73'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +1.76%  '
$ws.Range("E9").Value = '  -0.61%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0616'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.60%  '
$ws.Range("E11").Value = '  +1.29%  '
$ws.Range("D12").Value = '1.874.33'
$ws.Range("E12").Value = '  +0.48%  '
$ws.Range("D13").Value = '1.640.53'
$ws.Range("E13").Value = '  +0.19%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.578'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +4.59%  '
$ws.Range("E15").Value = '  +1.36%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '66.10'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +1.36%  '
$ws.Range("D17").Value = '27.989.21'
$ws.Range("E17").Value = '  +1.35%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '233.60'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +1.29%  '
$ws.Range("D19").Value = '0.0₃0726'
$ws.Range("E19").Value = '  +0.80%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.63'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.63%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.00'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.04%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '10.83'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +2.18%  '
$ws.Range("E23").Value = '  +0.23%  '
$ws.Range("E24").Value = '  -1.20%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '151.65'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.57%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '6.97'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +1.10%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '15.77'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +1.44%  '
$ws.Range("E28").Value = '  +0.12%  '
$ws.Range("E29").Value = '  -0.01%  '
$ws.Range("E30").Value = '  +0.77%  '
$ws.Range("E31").Value = '  +0.26%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.35'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +2.02%  '
$ws.Range("E33").Value = '  +0.23%  '
$ws.Range("D34").Value = '1.400.57'
$ws.Range("E34").Value = '  -5.50%  '
$ws.Range("E35").Value = '  +2.53%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.36'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.70%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.893'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +1.56%  '
$ws.Range("E38").Value = '  +0.84%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.558'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.48%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.918'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -1.96%  '
$ws.Range("E41").Value = '  -0.79%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.87'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +7.25%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '66.32'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -2.39%  '
$ws.Range("E45").Value = '  +2.17%  '
$ws.Range("E46").Value = '  +0.09%  '
$ws.Range("D47").Value = '1.783.05'
$ws.Range("E47").Value = '  +0.53%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '88.17'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.54%  '
$ws.Range("E49").Value = '  +1.35%  '
$ws.Range("E50").Value = '  +0.33%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '7.63'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.30%  '
